# Logic tree input file updated
# Insert a new row at row 7 (pushing existing rows 7-11 down to 8-12),
# populate it with the "Possible_Problem" branch content, set its row
# height, copy the wrap-text style used by the other Possible_Problem
# cells in column C, and move the active selection to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7; existing rows 7-11 shift to 8-12.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7.
$ws.Range("A7").Value = "Problem:Does the power steering system works normally at times,  or does it alternate between difficult and easy when the steering wheel is turned? (Please answer as:  Alternate between difficult & easy, Works normally at times"
$ws.Range("B7").Value = "Possible_Problem"
$ws.Range("C7").Value = "Possible_Problem:50% Steering Gear / Module Assembly`n25% Steering Intermediate Shaft`n15% Fuse`n10% Power Steering Module Wiring Damage"

# Match formatting used elsewhere for these "Possible_Problem" answer rows.
$ws.Range("C7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 259.2

# Update the saved selection/view to match the authored state.
$ws.Range("A4").Select()
